$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C110").Value = 7310
$ws.Range("C111:C252").Value = 7293
